# Populate the pitch-by-pitch "Catcher's View" (Pitch / Choice / Result) table
# and the related Exit Velo / Launch Angle / Result summary cells for each
# at-bat block in the postgame hitter report. Also normalizes the "Pitch Mix:"
# ordering for each at-bat to start with the most-thrown-this-AB pitch type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- At-bat starting row 10 (Inning 1) ---
$ws.Range("F10").Value = "FB"
$ws.Range("G10").Value = "Take"
$ws.Range("H10").Value = "Ball"
$ws.Range("M10").Value = "93.8 MPH"
$ws.Range("F11").Value = "FB"
$ws.Range("G11").Value = "Swing"
$ws.Range("H11").Value = "In Play"
$ws.Range("M12").Value = "39.82°"
$ws.Range("J17").Value = "CH,CB,FB"

# --- At-bat starting row 19 (Inning 3) ---
$ws.Range("F19").Value = "CB"
$ws.Range("G19").Value = "Take"
$ws.Range("H19").Value = "Ball"
$ws.Range("F20").Value = "FB"
$ws.Range("G20").Value = "Take"
$ws.Range("H20").Value = "Ball"
$ws.Range("F21").Value = "FB"
$ws.Range("G21").Value = "Take"
$ws.Range("H21").Value = "Ball"
$ws.Range("M21").Value = $null
$ws.Range("F22").Value = "CH"
$ws.Range("G22").Value = "Take"
$ws.Range("H22").Value = "Strike"
$ws.Range("F23").Value = "FB"
$ws.Range("G23").Value = "Take"
$ws.Range("H23").Value = "Strike"
$ws.Range("F24").Value = "FB"
$ws.Range("G24").Value = "Take"
$ws.Range("H24").Value = "Ball"
$ws.Range("M24").Value = "Walk"
$ws.Range("J26").Value = "CH,CB,FB"

# --- At-bat starting row 28 ---
$ws.Range("F28").Value = "CH"
$ws.Range("G28").Value = "Swing"
$ws.Range("H28").Value = "Foul"
$ws.Range("F29").Value = "CH"
$ws.Range("G29").Value = "Take"
$ws.Range("H29").Value = "Ball"
$ws.Range("F30").Value = "CH"
$ws.Range("G30").Value = "Take"
$ws.Range("H30").Value = "Ball"
$ws.Range("M30").Value = $null
$ws.Range("F31").Value = "FB"
$ws.Range("G31").Value = "Take"
$ws.Range("H31").Value = "Ball"
$ws.Range("F32").Value = "CH"
$ws.Range("G32").Value = "Swing"
$ws.Range("H32").Value = "Foul"
$ws.Range("F33").Value = "FB"
$ws.Range("G33").Value = "Take"
$ws.Range("H33").Value = "Ball"
$ws.Range("M33").Value = "Walk"
$ws.Range("J35").Value = "CH,CB,FB,SL"

# --- At-bat starting row 37 ---
$ws.Range("F37").Value = "CH"
$ws.Range("G37").Value = "Take"
$ws.Range("H37").Value = "Ball"
$ws.Range("F38").Value = "CH"
$ws.Range("G38").Value = "Swing"
$ws.Range("H38").Value = "Strike"
$ws.Range("F39").Value = "CH"
$ws.Range("G39").Value = "Take"
$ws.Range("H39").Value = "Ball"
$ws.Range("M39").Value = $null
$ws.Range("F40").Value = "CH"
$ws.Range("G40").Value = "Swing"
$ws.Range("H40").Value = "Foul"
$ws.Range("F41").Value = "CH"
$ws.Range("G41").Value = "Swing"
$ws.Range("H41").Value = "Foul"
$ws.Range("F42").Value = "CH"
$ws.Range("G42").Value = "Swing"
$ws.Range("H42").Value = "Foul"
$ws.Range("M42").Value = "Walk"
$ws.Range("F43").Value = "CB"
$ws.Range("G43").Value = "Take"
$ws.Range("H43").Value = "Ball"
$ws.Range("F44").Value = "CH"
$ws.Range("G44").Value = "Take"
$ws.Range("H44").Value = "Ball"
$ws.Range("J44").Value = "CH,CB,FB,SL"

# --- At-bat starting row 46 ---
$ws.Range("F46").Value = "FB"
$ws.Range("G46").Value = "Take"
$ws.Range("H46").Value = "Ball"
$ws.Range("M46").Value = "101.29 MPH"
$ws.Range("F47").Value = "CH"
$ws.Range("G47").Value = "Swing"
$ws.Range("H47").Value = "In Play"
$ws.Range("M48").Value = "11.83°"
$ws.Range("J53").Value = "CH,FB,SL"

# --- At-bat starting row 61 ---
$ws.Range("F61").Value = "CH"
$ws.Range("G61").Value = "Take"
$ws.Range("H61").Value = "Ball"
$ws.Range("M61").Value = "55.63 MPH"
$ws.Range("F62").Value = "CH"
$ws.Range("G62").Value = "Swing"
$ws.Range("H62").Value = "In Play"
$ws.Range("M63").Value = "-43.57°"
$ws.Range("J68").Value = "CH,FB,SL"
